# Improved canal ponding: distributed water throughout cell
# Add two new pointer rows (template_output_raster / initial_zeta_pickle) below
# the existing file pointer table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "template_output_raster"
$ws.Range("A11").Value = "initial_zeta_pickle"
$ws.Range("B11").Value = "data/new_area/best_initial_zeta.p"
$ws.Range("B12").Value = "data/dtm_depth_padded.tif"

$ws.Range("B17").Select()
